$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I ("I0") and J ("IF") ---
# Copy formatting from existing header cell H1 so the new header cells
# reuse the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Column I (rows 2-8): constant value 1 ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# --- Column J (rows 2-8): mirrors column H values ---
for ($r = 2; $r -le 8; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hVal
}
